# duplicate psx folder deleted
#
# The QSE sheet ended up with a handful of ticker symbols duplicated
# (CBQK, ABQK, DOHI, BEMA). Remove those duplicate rows.

$wb = $excel.ActiveWorkbook
$qse = $wb.Worksheets.Item("QSE")

$duplicateSymbols = @("CBQK", "ABQK", "DOHI", "BEMA")

# Find the row number for every duplicate symbol first ...
$usedRange = $qse.UsedRange
$rowCount = $usedRange.Rows.Count

$rowsToDelete = @()
for ($r = 1; $r -le $rowCount; $r++) {
    $val = $qse.Cells.Item($r, 1).Value()
    if ($duplicateSymbols -contains $val) {
        $rowsToDelete += $r
    }
}

# ... then delete them starting from the bottom so the row numbers found
# above stay valid while we work our way back up the sheet.
$rowsToDelete = $rowsToDelete | Sort-Object -Descending
foreach ($r in $rowsToDelete) {
    $qse.Rows.Item($r).Delete()
}

# Make QSE the active sheet/tab (it was CUSTUM before) and restore the
# view/selection state recorded for QSE after the edit.
$qse.Activate()
$qse.Range("A45").EntireRow.Select()
$excel.ActiveWindow.ScrollRow = 37

$wb.Save()
